$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("B2").Value = "C:\Users\W8DE5P2\OneDrive-Deere&Co\OneDrive - Deere & Co\Desktop\Proveedores\CLIIENTES JOHN DEERE\JD SARAN\t134.pdf"
$ws.Range("D2").Value = "99999TCD00"
$ws.Range("I2").Value = 67.855

# --- Duplicate row 2's formatting (borders/font/alignment) down into rows 3 and 4 ---
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A3:I3").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:I4").PasteSpecial(-4122) | Out-Null

# Columns C, F and G hold digit-only / date-look-alike text; force text format
# so Excel does not silently reinterpret them as numbers or dates.
$ws.Range("C3:C4").NumberFormat = "@"
$ws.Range("F3:F4").NumberFormat = "@"
$ws.Range("G3:G4").NumberFormat = "@"

# --- New row 3 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "C:\Users\W8DE5P2\OneDrive-Deere&Co\OneDrive - Deere & Co\Desktop\Proveedores\CLIIENTES JOHN DEERE\JD SARAN\t54.pdf"
$ws.Range("C3").Value = "4501168528"
$ws.Range("D3").Value = "99999TCD00"
$ws.Range("E3").Value = "R104907"
$ws.Range("F3").Value = "250"
$ws.Range("G3").Value = "07.12.2021"
$ws.Range("I3").Value = 79.55500000000001

# --- New row 4 ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "C:\Users\W8DE5P2\OneDrive-Deere&Co\OneDrive - Deere & Co\Desktop\Proveedores\CLIIENTES JOHN DEERE\JD SARAN\t66.pdf"
$ws.Range("C4").Value = "4501168528"
$ws.Range("D4").Value = "99999TCD00"
$ws.Range("E4").Value = "R104907"
$ws.Range("F4").Value = "250"
$ws.Range("G4").Value = "07.12.2021"
$ws.Range("I4").Value = 79.55500000000001
